# Add a final slide to the deck ("Title and Content" layout, same as the
# other content slides in this deck), with the closing "optional extra
# slides" text that was added to the author's deck.

$p = $ppt.ActivePresentation

# Layout index 2 = "Title and Content" (same layout slide 5 / slideLayout2.xml
# already uses in this deck).
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder (shape 1)
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Feel free to add one or two more slides to discuss whatever else you" + [char]0x2019 + "d like"

# Body / content placeholder (shape 2) -- build it as two runs, "(" then
# "not required)", matching how the text was typed up in the original edit.
$body = $newSlide.Shapes.Item(2)
$bodyRange = $body.TextFrame.TextRange
$bodyRange.Text = "("
$null = $bodyRange.InsertAfter("not required)")
